$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.158.17"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "1.857.43"
$ws.Range("E3").Value = "  -3.41%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'234.26"
$ws.Range("E5").Value = "  -2.98%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.4664"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("D8").Value = "'0.2821"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").Value = "'0.06564"
$ws.Range("E9").Value = "  -3.13%  "
$ws.Range("D10").Value = "'19.99"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").Value = "'0.07832"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'96.88"
$ws.Range("E12").Value = "  -7.29%  "
$ws.Range("D13").Value = "1.867.59"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "'5.115"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").Value = "'0.6654"
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").Value = "'282.40"
$ws.Range("E16").Value = "  -4.25%  "
$ws.Range("D17").Value = "30.195.49"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "'5.438"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "'12.59"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").Value = "2.111.61"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("D22").Value = "'0.000007233"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'6.143"
$ws.Range("E24").Value = "  -3.83%  "
$ws.Range("D25").Value = "'9.333"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").Value = "'167.65"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").Value = "'18.93"
$ws.Range("E27").Value = "  -4.40%  "
$ws.Range("D28").Value = "'1.916"
$ws.Range("E28").Value = "  -9.45%  "
$ws.Range("D29").Value = "'1.339"
$ws.Range("E29").Value = "  -3.52%  "
$ws.Range("D30").Value = "'0.09572"
$ws.Range("E30").Value = "  -4.75%  "
$ws.Range("D31").Value = "'4.408"
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("D32").Value = "'1.469"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").Value = "'4.100"
$ws.Range("E33").Value = "  -4.96%  "
$ws.Range("D34").Value = "'0.04675"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").Value = "'0.7004"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").Value = "'1.098"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").Value = "'0.9999"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'2.707"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "'0.01853"
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("D40").Value = "'6.355"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").Value = "'2.511"
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("D42").Value = "'72.03"
$ws.Range("E42").Value = "  -4.29%  "
$ws.Range("D43").Value = "'0.8531"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").Value = "'1.931"
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'0.4162"
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").Value = "'103.67"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").Value = "'991.53"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").Value = "'7.214"
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("D50").Value = "'9.154"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'33.93"
$ws.Range("E51").Value = "  -3.08%  "
